# Config.xlsx update
# - Adds a new "DiccionarioNavegacionSiesa" parameter row to the Config sheet
# - Updates TamanioDeLote value from 200 to 100
# - Widens column B on Config to fit the new long dictionary value
# - Makes "Config" the active/selected sheet (was "Notificaciones")
# - Adjusts the header row height on "Notificaciones" (wrap reflow)

$wb = $excel.ActiveWorkbook

$wsConfig = $wb.Worksheets.Item("Config")
$wsNotif  = $wb.Worksheets.Item("Notificaciones")

# --- Config sheet edits ---------------------------------------------------

# TamanioDeLote: 200 -> 100
$wsConfig.Range("B8").Value = 100

# Widen column B so the long dictionary string is legible
$wsConfig.Columns.Item(2).ColumnWidth = 94.5546875

# New row 12: DiccionarioNavegacionSiesa parameter
$wsConfig.Range("A12").Value = "DiccionarioNavegacionSiesa"
$wsConfig.Range("B12").Value = "{'Financiero': 'f','Contabilidad General': 'g','Contab': 'n','Auditoria de Documentos': 'r','Consulta por Numero': 'n'}"

# Match formatting used for the other long-text cells: wrap + vertical center
$wsConfig.Range("B12").WrapText = $true
$wsConfig.Range("B12").VerticalAlignment = -4108

# Row height for the new row
$wsConfig.Rows.Item(12).RowHeight = 17.4

# --- Notificaciones sheet edits -------------------------------------------

# Header row shrinks after wrap reflow
$wsNotif.Rows.Item(1).RowHeight = 26.4

# --- Active sheet / selection ----------------------------------------------
# The workbook now opens on Config rather than Notificaciones
$null = $wsConfig.Activate()
$null = $wsConfig.Range("B17").Select()
